# Weekly refresh of the Espinaca price table: the source feed re-pulled the
# week's records, which landed in different rows than before. Re-map the
# per-row fields (Fecha, Calidad, Volumen, Precio mínimo/máximo/promedio,
# Origen, Precio $/Kg) from their old row to their new row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map: new row number -> old row number (where the record used to live).
$rowMap = @{
    2  = 19
    3  = 22
    4  = 3
    5  = 17
    6  = 21
    7  = 8
    8  = 12
    9  = 2
    10 = 20
    11 = 9
    12 = 14
    13 = 15
    14 = 16
    15 = 23
    16 = 7
    17 = 18
    18 = 13
    19 = 11
    20 = 6
    21 = 5
    22 = 4
    23 = 10
}

# Snapshot every affected column (D, I, J, K, L, M, O, P) from every data row
# *before* writing anything back, so the permutation below never reads a
# value that's already been overwritten.
$snapshot = @{}
foreach ($r in 2..23) {
    $snapshot[$r] = @{
        D = $ws.Cells.Item($r, 4).Value2
        I = $ws.Cells.Item($r, 9).Value2
        J = $ws.Cells.Item($r, 10).Value2
        K = $ws.Cells.Item($r, 11).Value2
        L = $ws.Cells.Item($r, 12).Value2
        M = $ws.Cells.Item($r, 13).Value2
        O = $ws.Cells.Item($r, 15).Value2
        P = $ws.Cells.Item($r, 16).Value2
    }
}

foreach ($newRow in $rowMap.Keys) {
    $oldRow = $rowMap[$newRow]
    $src = $snapshot[$oldRow]

    $ws.Cells.Item($newRow, 4).Value2 = $src.D    # Fecha
    $ws.Cells.Item($newRow, 9).Value2 = $src.I    # Calidad
    $ws.Cells.Item($newRow, 10).Value2 = $src.J   # Volumen
    $ws.Cells.Item($newRow, 11).Value2 = $src.K   # Precio minimo
    $ws.Cells.Item($newRow, 12).Value2 = $src.L   # Precio maximo
    $ws.Cells.Item($newRow, 13).Value2 = $src.M   # Precio promedio ponderado
    $ws.Cells.Item($newRow, 15).Value2 = $src.O   # Origen
    $ws.Cells.Item($newRow, 16).Value2 = $src.P   # Precio $/Kg
}
